$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '246.52'
$ws.Range("E2").Value = '1.09%'

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '29.45'
$ws.Range("E3").Value = '-2.58%'

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '5.156'
$ws.Range("E4").Value = '0.14%'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.12%'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.57%'

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '3.188'
$ws.Range("E7").Value = '5.24%'

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8518'
$ws.Range("E8").Value = '0.69%'

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8663'
$ws.Range("E9").Value = '1.08%'

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1371'
$ws.Range("E10").Value = '1.44%'

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.07069'
$ws.Range("E11").Value = '2.23%'

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03264'
$ws.Range("E12").Value = '13.06%'

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09366'
$ws.Range("E13").Value = '-0.09%'

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001528'
$ws.Range("E14").Value = '-0.47%'

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = '0.0005972'
$ws.Range("E15").Value = '-0.20%'

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005961'
$ws.Range("E16").Value = '-4.10%'

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = '3.486'
$ws.Range("E17").Value = '-0.65%'

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = '2.212'
$ws.Range("E18").Value = '-0.95%'

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3200'
$ws.Range("E19").Value = '1.61%'

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03353'
$ws.Range("E20").Value = '0.50%'

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '3.147'
$ws.Range("E22").Value = '-13.10%'

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04129'
$ws.Range("E23").Value = '-1.07%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.85%'

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001225'
$ws.Range("E25").Value = '1.26%'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-6.77%'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.48%'

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001443'
$ws.Range("E28").Value = '3.84%'

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03743'
$ws.Range("E40").Value = '-0.68%'

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005772'
$ws.Range("E41").Value = '8.38%'

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1072'
$ws.Range("E42").Value = '1.31%'

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002406'
$ws.Range("E43").Value = '5.20%'

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009154'
$ws.Range("E44").Value = '-1.38%'

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005281'
$ws.Range("E45").Value = '3.58%'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.05%'

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05792'
$ws.Range("E47").Value = '-42.02%'

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002173'
$ws.Range("E48").Value = '-21.80%'

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002097'
$ws.Range("E49").Value = '-0.05%'

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001997'
$ws.Range("E50").Value = '-0.05%'
